# Customer.xlsx -- "added Customer Creation Flow in app"
#
# Sheet1 row 5 (customerVillage) held the shared string "Belakavadi S.O"
# (shared by row 4 / customerPostal too). The edit introduces a distinct
# value for the village field: a brand-new string "Belakavadi", while the
# postal field keeps "Belakavadi S.O". The active selection also moves to
# the edited cell (B5) with the sheet scrolled back to the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New, separate shared-string value for customerVillage.
$ws.Range("B5").Value = "Belakavadi"

# Scroll back to the top and land the selection on the edited cell.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B5").Select()
